$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct values: C column holds "Npaginas" per user; D column ("Mpag") should hold
# the MAX of Npaginas within each merged "equipo" group, not the SUM.

# Group "DP VESTIR" (rows 14-16)
$ws.Range("C14").Value = 35
$ws.Range("D14").Value = 35

$ws.Range("C15").Value = 16
$ws.Range("D15").Value = 35

$ws.Range("D16").Value = 35

# Group "EC" (rows 19-21)
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0

$ws.Range("D20").Value = 0

$ws.Range("D21").Value = 0
